$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.929.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8266"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3209"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06995"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08031"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.913.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7445"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.178"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.932.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.892"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007733"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.154.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.892"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1571"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +22.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.081"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.376"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.519"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.233"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05649"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.059"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.268"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7287"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.50%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01901"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.767"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4389"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.927"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8454"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("E44").Value = "  +0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.883"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.564"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.665"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "984.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.052.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
